$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.414.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.328.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.47"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.60"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0917"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.680.93"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.325.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.296.35"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -10.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.75"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.60"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +9.48%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.95"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.60"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.24"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0886"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.77"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.69%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0362"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.79"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.34"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +13.11%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.63"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.235"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.30"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.13%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.01"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.661.61"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.73"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.33"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.77%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.93"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.43%  "
